# Updates cryptos list prices and 1h trading volume percentages
# (scheduled GitHub Actions refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.934.06'
$ws.Range('E2').Value = '  +0.04%  '

$ws.Range('D3').Value = '3.182.03'
$ws.Range('E3').Value = '  -0.59%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').Value = '604.57'
$ws.Range('E5').Value = '  +1.53%  '

$ws.Range('D6').Value = '154.51'
$ws.Range('E6').Value = '  +0.43%  '

$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').Value = '3.181.37'
$ws.Range('E8').Value = '  -0.57%  '

$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  +1.94%  '

$ws.Range('E10').Value = '  -1.44%  '

$ws.Range('D11').Value = '5.65'
$ws.Range('E11').Value = '  -6.95%  '

$ws.Range('D12').Value = '0.513'
$ws.Range('E12').Value = '  -0.07%  '

$ws.Range('D13').Value = '0.0000264'
$ws.Range('E13').Value = '  -2.71%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '38.20'
$ws.Range('E14').Value = '  -2.04%  '

$ws.Range('D15').Value = '3.698.52'
$ws.Range('E15').Value = '  -0.64%  '

$ws.Range('D16').Value = '65.955.42'
$ws.Range('E16').Value = '  +0.07%  '

$ws.Range('D17').Value = '7.38'
$ws.Range('E17').Value = '  -0.53%  '

$ws.Range('D18').Value = '3.176.44'
$ws.Range('E18').Value = '  -0.74%  '

$ws.Range('E19').Value = '  +1.01%  '

$ws.Range('D20').Value = '507.06'
$ws.Range('E20').Value = '  -0.57%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.30'
$ws.Range('E21').Value = '  -0.11%  '

$ws.Range('D22').Value = '0.728'
$ws.Range('E22').Value = '  -1.67%  '

$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('E24').Value = '  -3.29%  '

$ws.Range('D25').Value = '84.28'
$ws.Range('E25').Value = '  -0.75%  '

$ws.Range('E26').Value = '  +0.20%  '

$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').Value = '9.16'
$ws.Range('E28').Value = '  -2.11%  '

$ws.Range('E29').Value = '  +4.91%  '

$ws.Range('D30').Value = '3.01'
$ws.Range('E30').Value = '  +4.46%  '

$ws.Range('E31').Value = '  +3.98%  '

$ws.Range('D32').Value = '27.92'
$ws.Range('E32').Value = '  -1.52%  '

$ws.Range('E33').Value = '  +0.14%  '

$ws.Range('E34').Value = '  -4.64%  '

$ws.Range('D35').Value = '6.48'
$ws.Range('E35').Value = '  -1.24%  '

$ws.Range('D36').Value = '513.75'
$ws.Range('E36').Value = '  +6.48%  '

$ws.Range('D37').Value = '55.18'
$ws.Range('E37').Value = '  +0.28%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0880'
$ws.Range('E38').Value = '  -2.87%  '

$ws.Range('D39').Value = '0.0416'
$ws.Range('E39').Value = '  -0.83%  '

$ws.Range('E40').Value = '  +5.13%  '

$ws.Range('D41').Value = '0.0₃0693'
$ws.Range('E41').Value = '  +6.91%  '

$ws.Range('D42').Value = '8.72'
$ws.Range('E42').Value = '  -1.26%  '

$ws.Range('D43').Value = '2.86'
$ws.Range('E43').Value = '  -2.20%  '

$ws.Range('E44').Value = '  -0.20%  '

$ws.Range('D45').Value = '2.48'
$ws.Range('E45').Value = '  +2.09%  '

$ws.Range('D46').Value = '2.832.27'
$ws.Range('E46').Value = '  -3.55%  '

$ws.Range('D47').Value = '27.81'
$ws.Range('E47').Value = '  -1.82%  '

$ws.Range('D49').Value = '2.37'
$ws.Range('E49').Value = '  +3.53%  '

$ws.Range('E50').Value = '  +0.26%  '

$ws.Range('D51').Value = '2.61'
$ws.Range('E51').Value = '  +2.13%  '
